$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.448.86'
$ws.Range('E2').Value = '  -2.11%  '
$ws.Range('D3').Value = '3.460.97'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '581.84'
$ws.Range('E5').Value = '  +5.19%  '
$ws.Range('D6').Value = '176.41'
$ws.Range('E6').Value = '  -5.69%  '
$ws.Range('E7').Value = '  +3.56%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('E10').Value = '  +5.37%  '
$ws.Range('D11').Value = '55.50'
$ws.Range('E11').Value = '  +1.82%  '
$ws.Range('E12').Value = '  +2.34%  '
$ws.Range('E13').Value = '  -1.33%  '
$ws.Range('D14').Value = '4.015.13'
$ws.Range('B15').Value = 'TRON'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D15').Value = '0.121'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.463.84'
$ws.Range('E16').Value = '  -2.12%  '
$ws.Range('D17').Value = '18.20'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').Value = '12.00'
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('D19').Value = '65.435.75'
$ws.Range('E19').Value = '  -2.15%  '
$ws.Range('E20').Value = '  +1.17%  '
$ws.Range('D21').Value = '409.88'
$ws.Range('E21').Value = '  -5.08%  '
$ws.Range('D22').Value = '4.28'
$ws.Range('E22').Value = '  +4.69%  '
$ws.Range('D23').Value = '4.49'
$ws.Range('E23').Value = '  +9.82%  '
$ws.Range('D24').Value = '84.45'
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('D25').Value = '13.41'
$ws.Range('E25').Value = '  +10.36%  '
$ws.Range('D26').Value = '10.99'
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('D27').Value = '2.85'
$ws.Range('E27').Value = '  -1.25%  '
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('D29').Value = '30.15'
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('D30').Value = '6.69'
$ws.Range('E30').Value = '  +1.66%  '
$ws.Range('D31').Value = '11.68'
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('D32').Value = '591.19'
$ws.Range('E32').Value = '  -8.75%  '
$ws.Range('E33').Value = '  -1.54%  '
$ws.Range('D34').Value = '60.62'
$ws.Range('E34').Value = '  +1.74%  '
$ws.Range('D35').Value = '0.152'
$ws.Range('E35').Value = '  +0.70%  '
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').Value = '0.0₃0787'
$ws.Range('E37').Value = '  -4.39%  '
$ws.Range('D38').Value = '36.81'
$ws.Range('E38').Value = '  -4.34%  '
$ws.Range('D39').Value = '3.58'
$ws.Range('E39').Value = '  +6.09%  '
$ws.Range('E40').Value = '  -2.05%  '
$ws.Range('D41').Value = '3.216.84'
$ws.Range('E41').Value = '  +5.65%  '
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('E43').Value = '  +2.74%  '
$ws.Range('E44').Value = '  -1.35%  '
$ws.Range('D45').Value = '2.53'
$ws.Range('E45').Value = '  -5.46%  '
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('E47').Value = '  +1.15%  '
$ws.Range('E48').Value = '  -5.42%  '
$ws.Range('D49').Value = '8.59'
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('D50').Value = '138.48'
$ws.Range('E50').Value = '  -1.74%  '
$ws.Range('E51').Value = '  -2.62%  '
